# Auto-generated: updates market-price derived columns (H-N) across all 8 job sheets
# per the scheduled-runner data refresh described in the commit message.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3859.8
$ws.Range("I62").Value = 3999
$ws.Range("J62").Value = 3825
$ws.Range("K62").Value = 3999
$ws.Range("L62").Value = 3825
$ws.Range("M62").Value = -3375
$ws.Range("N62").Value = -5073
$ws.Range("H65").Value = 3859.8
$ws.Range("I65").Value = 3999
$ws.Range("J65").Value = 3825
$ws.Range("K65").Value = 19995
$ws.Range("L65").Value = 19125
$ws.Range("M65").Value = -16875
$ws.Range("N65").Value = -25365
$ws.Range("H98").Value = 1197.8823
$ws.Range("I98").Value = 1085.25
$ws.Range("K98").Value = 1085.25
$ws.Range("M98").Value = 412.75
$ws.Range("H122").Value = 1197.8823
$ws.Range("I122").Value = 1085.25
$ws.Range("K122").Value = 3255.75
$ws.Range("M122").Value = -805.75
$ws.Range("H127").Value = 1290.8667
$ws.Range("J127").Value = 1338.5834
$ws.Range("L127").Value = 4015.7502
$ws.Range("N127").Value = -13935.7502
$ws.Range("H132").Value = 4008.6296
$ws.Range("I132").Value = 3439.6191
$ws.Range("J132").Value = 6000.1665
$ws.Range("K132").Value = 10318.8573
$ws.Range("L132").Value = 18000.4995
$ws.Range("M132").Value = -7788.8573
$ws.Range("N132").Value = -23060.4995
$ws.Range("H137").Value = 3523.3684
$ws.Range("I137").Value = 3523.3684
$ws.Range("K137").Value = 10570.1052
$ws.Range("M137").Value = -8020.1052
$ws.Range("H138").Value = 166832.47
$ws.Range("I138").Value = 1799.32
$ws.Range("J138").Value = 267462.44
$ws.Range("K138").Value = 5397.96
$ws.Range("L138").Value = 802387.3200000001
$ws.Range("M138").Value = -257.96
$ws.Range("N138").Value = -812667.3200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 438216.03
$ws.Range("I32").Value = 490634.8
$ws.Range("J32").Value = 18866
$ws.Range("K32").Value = 490634.8
$ws.Range("L32").Value = 18866
$ws.Range("M32").Value = -490347.8
$ws.Range("N32").Value = -19440
$ws.Range("H43").Value = 14377
$ws.Range("J43").Value = 14377
$ws.Range("L43").Value = 14377
$ws.Range("N43").Value = -15003
$ws.Range("H63").Value = 5704.5625
$ws.Range("I63").Value = 3726.7144
$ws.Range("J63").Value = 7242.8887
$ws.Range("K63").Value = 3726.7144
$ws.Range("L63").Value = 7242.8887
$ws.Range("M63").Value = -3040.7144
$ws.Range("N63").Value = -8614.8887
$ws.Range("H66").Value = 5704.5625
$ws.Range("I66").Value = 3726.7144
$ws.Range("J66").Value = 7242.8887
$ws.Range("K66").Value = 18633.572
$ws.Range("L66").Value = 36214.4435
$ws.Range("M66").Value = -15201.572
$ws.Range("N66").Value = -43078.4435
$ws.Range("H122").Value = 2009
$ws.Range("I122").Value = 1945.7778
$ws.Range("J122").Value = 2122.8
$ws.Range("K122").Value = 5837.3334
$ws.Range("L122").Value = 6368.400000000001
$ws.Range("M122").Value = -3387.3334
$ws.Range("N122").Value = -11268.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2880.8147
$ws.Range("I134").Value = 2275.4119
$ws.Range("K134").Value = 6826.2357
$ws.Range("M134").Value = -4291.2357

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 21752
$ws.Range("J2").Value = 43000
$ws.Range("L2").Value = 43000
$ws.Range("N2").Value = -43226
$ws.Range("H22").Value = 469.16666
$ws.Range("I22").Value = 278.75
$ws.Range("J22").Value = 850
$ws.Range("K22").Value = 278.75
$ws.Range("L22").Value = 850
$ws.Range("M22").Value = 71.25
$ws.Range("N22").Value = -1550
$ws.Range("H31").Value = 2769.868
$ws.Range("I31").Value = 871.2
$ws.Range("J31").Value = 6461.722
$ws.Range("K31").Value = 871.2
$ws.Range("L31").Value = 6461.722
$ws.Range("M31").Value = -576.2
$ws.Range("N31").Value = -7051.722
$ws.Range("H34").Value = 2769.868
$ws.Range("I34").Value = 871.2
$ws.Range("J34").Value = 6461.722
$ws.Range("K34").Value = 871.2
$ws.Range("L34").Value = 6461.722
$ws.Range("M34").Value = -669.2
$ws.Range("N34").Value = -6865.722
$ws.Range("H107").Value = 2717861
$ws.Range("I107").Value = 4167052.8
$ws.Range("J107").Value = 626.5
$ws.Range("K107").Value = 4167052.8
$ws.Range("L107").Value = 626.5
$ws.Range("M107").Value = -4165132.8
$ws.Range("N107").Value = -4466.5
$ws.Range("H122").Value = 1892.7273
$ws.Range("I122").Value = 1478
$ws.Range("J122").Value = 1984.8889
$ws.Range("K122").Value = 4434
$ws.Range("L122").Value = 5954.6667
$ws.Range("N122").Value = -10854.6667
$ws.Range("M122").Value = -1984
$ws.Range("H132").Value = 16669125
$ws.Range("I132").Value = 945
$ws.Range("K132").Value = 2835
$ws.Range("M132").Value = -305

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 945.25806
$ws.Range("J131").Value = 1043.4445
$ws.Range("L131").Value = 3130.3335
$ws.Range("N131").Value = -13210.3335
$ws.Range("H136").Value = 1772.5
$ws.Range("I136").Value = 1454.2858
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 4362.857400000001
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = 737.1425999999992
$ws.Range("N136").Value = -22200
$ws.Range("H138").Value = 2503.2
$ws.Range("I138").Value = 1630
$ws.Range("J138").Value = 2677.84
$ws.Range("K138").Value = 4890
$ws.Range("L138").Value = 8033.52
$ws.Range("M138").Value = 250
$ws.Range("N138").Value = -18313.52
$ws.Range("H139").Value = 3437.875
$ws.Range("I139").Value = 1334
$ws.Range("K139").Value = 4002
$ws.Range("M139").Value = 1138
$ws.Range("H140").Value = 1764
$ws.Range("I140").Value = 1159.625
$ws.Range("K140").Value = 3478.875
$ws.Range("M140").Value = 1701.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 79042
$ws.Range("J26").Value = 79042
$ws.Range("L26").Value = 79042
$ws.Range("N26").Value = -79602
$ws.Range("H28").Value = 29999
$ws.Range("J28").Value = 29999
$ws.Range("L28").Value = 29999
$ws.Range("N28").Value = -30383
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = ""
$ws.Range("H41").Value = 2357.6667
$ws.Range("I41").Value = 3124
$ws.Range("K41").Value = 3124
$ws.Range("M41").Value = -2769
$ws.Range("H50").Value = 79042
$ws.Range("J50").Value = 79042
$ws.Range("L50").Value = 79042
$ws.Range("N50").Value = -80038
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""
$ws.Range("H122").Value = 5118.385
$ws.Range("I122").Value = 4285.4287
$ws.Range("J122").Value = 5425.263
$ws.Range("K122").Value = 12856.2861
$ws.Range("L122").Value = 16275.789
$ws.Range("M122").Value = -10406.2861
$ws.Range("N122").Value = -21175.789
$ws.Range("H132").Value = 1893.2593
$ws.Range("I132").Value = 1576.4762
$ws.Range("J132").Value = 3002
$ws.Range("K132").Value = 4729.4286
$ws.Range("L132").Value = 9006
$ws.Range("M132").Value = -2199.4286
$ws.Range("N132").Value = -14066

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 10645.818
$ws.Range("I93").Value = 17593.5
$ws.Range("J93").Value = 2308.6
$ws.Range("K93").Value = 17593.5
$ws.Range("L93").Value = 2308.6
$ws.Range("M93").Value = -16345.5
$ws.Range("N93").Value = -4804.6
$ws.Range("H122").Value = 5310.3335
$ws.Range("I122").Value = 6158.8
$ws.Range("J122").Value = 4249.75
$ws.Range("K122").Value = 18476.4
$ws.Range("L122").Value = 12749.25
$ws.Range("M122").Value = -16026.4
$ws.Range("N122").Value = -17649.25
$ws.Range("H132").Value = 4166.8716
$ws.Range("I132").Value = 3653.238
$ws.Range("J132").Value = 4766.1113
$ws.Range("K132").Value = 10959.714
$ws.Range("L132").Value = 14298.3339
$ws.Range("M132").Value = -8429.714
$ws.Range("N132").Value = -19358.3339
$ws.Range("H136").Value = 5954263
$ws.Range("I136").Value = 2077.375
$ws.Range("K136").Value = 6232.125
$ws.Range("M136").Value = -3682.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 63179
$ws.Range("I62").Value = 3902
$ws.Range("J62").Value = 77998.25
$ws.Range("K62").Value = 3902
$ws.Range("L62").Value = 77998.25
$ws.Range("M62").Value = -3278
$ws.Range("N62").Value = -79246.25
$ws.Range("H65").Value = 63179
$ws.Range("I65").Value = 3902
$ws.Range("J65").Value = 77998.25
$ws.Range("K65").Value = 19510
$ws.Range("L65").Value = 389991.25
$ws.Range("M65").Value = -16390
$ws.Range("N65").Value = -396231.25
$ws.Range("H81").Value = 5636.9
$ws.Range("I81").Value = 7173.8
$ws.Range("J81").Value = 4100
$ws.Range("K81").Value = 14347.6
$ws.Range("L81").Value = 8200
$ws.Range("M81").Value = -13286.6
$ws.Range("N81").Value = -10322
$ws.Range("H84").Value = 5636.9
$ws.Range("I84").Value = 7173.8
$ws.Range("J84").Value = 4100
$ws.Range("K84").Value = 71738
$ws.Range("L84").Value = 41000
$ws.Range("M84").Value = -66434
$ws.Range("N84").Value = -51608
$ws.Range("H100").Value = 2229.6667
$ws.Range("I100").Value = 699
$ws.Range("J100").Value = 2995
$ws.Range("K100").Value = 1398
$ws.Range("L100").Value = 5990
$ws.Range("M100").Value = -857
$ws.Range("N100").Value = -7072
$ws.Range("H107").Value = 679.3333
$ws.Range("I107").Value = 739.2
$ws.Range("J107").Value = 380
$ws.Range("K107").Value = 2217.6
$ws.Range("L107").Value = 1140
$ws.Range("M107").Value = -297.6000000000004
$ws.Range("N107").Value = -4980
$ws.Range("H122").Value = 3633.3333
$ws.Range("I122").Value = 2950
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 8850
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900
$ws.Range("M122").Value = -6400
$ws.Range("H132").Value = 3088745
$ws.Range("I132").Value = 2584.6553
$ws.Range("J132").Value = 6668691
$ws.Range("K132").Value = 7753.965899999999
$ws.Range("L132").Value = 20006073
$ws.Range("M132").Value = -5223.965899999999
$ws.Range("N132").Value = -20011133
